# delete iteration option, added error handling for user inputs - Reem
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("C2:C4").NumberFormat = "@"

$ws.Range("A2").Value = "Itr 1"
$ws.Range("A3").Value = "Itr 2"
$ws.Range("A4").Value = "Itr 3"

$ws.Range("B2").Value = "14"
$ws.Range("B3").Value = "20"
$ws.Range("B4").Value = "20"

$ws.Range("C2").Value = "01/04/2021"
$ws.Range("C3").Value = "10/04/2021"
$ws.Range("C4").Value = "01/03/2021"

$ws.Range("D2").Value = "Reem-George-Shadi"
$ws.Range("D3").Value = "Reem-George"
$ws.Range("D4").Value = "Shadi-Wafic"

$ws.Range("B6").Select()
